# Changement d'ordre des couleurs
# Swap the "carreaux" and "cœur" labels in the "nom" column (D) of Feuil1.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Feuil1")

# Row 3: "AS de carreaux" -> "AS de cœur"
$ws.Range("D3").Value = "AS de cœur"
# Row 5: "As de cœur" -> "As de carreaux"
$ws.Range("D5").Value = "As de carreaux"
# Row 7: "2 de carreaux" -> "2 de cœur"
$ws.Range("D7").Value = "2 de cœur"
# Row 9: "2 de cœur" -> "2 de carreaux"
$ws.Range("D9").Value = "2 de carreaux"

# Update the active cell selection shown in the sheet view.
$ws.Range("D10").Select()
